$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-74 from 45171 to 45172
for ($row = 2; $row -le 74; $row++) {
    $ws.Cells.Item($row, 3).Value = 45172
}
